$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format before writing numeric-looking price strings
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '23.428.11'
$ws.Range('E2').Value = '  +1.11%  '

$ws.Range('D3').Value = '1.637.88'
$ws.Range('E3').Value = '  +2.32%  '

$ws.Range('D4').Value = '1.001'

$ws.Range('D5').Value = '1.001'
$ws.Range('E5').Value = '  +0.02%  '

$ws.Range('D6').Value = '304.98'
$ws.Range('E6').Value = '  +0.72%  '

$ws.Range('E7').Value = '  -1.04%  '

$ws.Range('D8').Value = '51.81'
$ws.Range('E8').Value = '  -0.19%  '

$ws.Range('D9').Value = '0.3622'
$ws.Range('E9').Value = '  +0.20%  '

$ws.Range('D10').Value = '1.256'
$ws.Range('E10').Value = '  -0.65%  '

$ws.Range('D11').Value = '0.08125'
$ws.Range('E11').Value = '  +0.20%  '

$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.05%  '

$ws.Range('D13').Value = '22.84'
$ws.Range('E13').Value = '  +0.47%  '

$ws.Range('D14').Value = '6.609'
$ws.Range('E14').Value = '  +0.47%  '

$ws.Range('D15').Value = '0.00001269'
$ws.Range('E15').Value = '  +2.27%  '

$ws.Range('D16').Value = '7.285'
$ws.Range('E16').Value = '  -1.45%  '

$ws.Range('D17').Value = '1.638.19'
$ws.Range('E17').Value = '  +2.52%  '

$ws.Range('D18').Value = '94.28'
$ws.Range('E18').Value = '  +0.47%  '

$ws.Range('D19').Value = '0.06910'
$ws.Range('E19').Value = '  +0.42%  '

$ws.Range('D20').Value = '18.12'
$ws.Range('E20').Value = '  +0.40%  '

$ws.Range('D21').Value = '6.509'
$ws.Range('E21').Value = '  -0.34%  '

$ws.Range('E22').Value = '  -0.01%  '

$ws.Range('D23').Value = '23.429.07'
$ws.Range('E23').Value = '  +1.13%  '

$ws.Range('D24').Value = '12.73'

$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').Value = '3.054'
$ws.Range('E25').Value = '  +2.49%  '

$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '2.421'
$ws.Range('E26').Value = '  +0.54%  '

$ws.Range('D27').Value = '21.21'
$ws.Range('E27').Value = '  +0.01%  '

$ws.Range('D28').Value = '151.61'
$ws.Range('E28').Value = '  +1.07%  '

$ws.Range('D29').Value = '5.320'
$ws.Range('E29').Value = '  +1.44%  '

$ws.Range('D30').Value = '136.02'

$ws.Range('D31').Value = '2.268'
$ws.Range('E31').Value = '  -4.17%  '

$ws.Range('D32').Value = '1.816.86'
$ws.Range('E32').Value = '  +2.10%  '

$ws.Range('D33').Value = '6.743'
$ws.Range('E33').Value = '  +0.04%  '

$ws.Range('D34').Value = '0.9547'
$ws.Range('E34').Value = '  -0.97%  '

$ws.Range('D35').Value = '0.02823'
$ws.Range('E35').Value = '  +4.13%  '

$ws.Range('D36').Value = '10.27'
$ws.Range('E36').Value = '  +0.56%  '

$ws.Range('D37').Value = '0.07264'
$ws.Range('E37').Value = '  -2.61%  '

$ws.Range('D38').Value = '0.2518'
$ws.Range('E38').Value = '  +0.33%  '

$ws.Range('D39').Value = '0.08785'
$ws.Range('E39').Value = '  -0.13%  '

$ws.Range('D40').Value = '6.076'
$ws.Range('E40').Value = '  +0.45%  '

$ws.Range('D41').Value = '1.374'
$ws.Range('E41').Value = '  +0.99%  '

$ws.Range('D42').Value = '0.7042'
$ws.Range('E42').Value = '  -0.65%  '

$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '12.44'
$ws.Range('E43').Value = '  -0.05%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '16.05'
$ws.Range('E44').Value = '  +3.36%  '

$ws.Range('D45').Value = '0.6505'
$ws.Range('E45').Value = '  -0.17%  '

$ws.Range('D46').Value = '2.324'
$ws.Range('E46').Value = '  +0.78%  '

$ws.Range('E47').Value = '  +0.09%  '

$ws.Range('D48').Value = '4.013'
$ws.Range('E48').Value = '  -0.07%  '

$ws.Range('D49').Value = '0.07981'
$ws.Range('E49').Value = '  +0.41%  '

$ws.Range('D50').Value = '128.28'
$ws.Range('E50').Value = '  -2.81%  '

$ws.Range('D51').Value = '1.203'
$ws.Range('E51').Value = '  +0.30%  '

# Restore default General formatting on column D (values remain text since they contain extra dots / were entered as text)
$ws.Range('D2:D51').ClearFormats()
